$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.607.22"
$ws.Range("E2").Value = "  +6.98%  "

$ws.Range("D3").Value = "2.638.30"
$ws.Range("E3").Value = "  +6.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.01"
$ws.Range("E5").Value = "  +12.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "584.39"
$ws.Range("E6").Value = "  +2.96%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +3.82%  "

$ws.Range("E9").Value = "  +11.96%  "

$ws.Range("D10").Value = "2.636.76"
$ws.Range("E10").Value = "  +6.92%  "

$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  +6.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.70"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").Value = "74.480.58"
$ws.Range("E14").Value = "  +7.01%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.120.91"
$ws.Range("E15").Value = "  +6.92%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("E16").Value = "  +3.07%  "

$ws.Range("E17").Value = "  +11.25%  "

$ws.Range("D18").Value = "2.621.35"
$ws.Range("E18").Value = "  +6.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.24"
$ws.Range("E19").Value = "  +28.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  +9.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.14"
$ws.Range("E21").Value = "  +6.67%  "

$ws.Range("E22").Value = "  +14.33%  "

$ws.Range("E23").Value = "  +4.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.17"
$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.09"
$ws.Range("E26").Value = "  +5.80%  "

$ws.Range("E27").Value = "  +6.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.31"
$ws.Range("E28").Value = "  +8.95%  "

$ws.Range("D29").Value = "2.757.47"
$ws.Range("E29").Value = "  +6.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("E30").Value = "  +2.64%  "

$ws.Range("E31").Value = "  +10.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "525.52"
$ws.Range("E32").Value = "  +17.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.40"
$ws.Range("E33").Value = "  +11.18%  "

$ws.Range("E34").Value = "  +4.55%  "

$ws.Range("E35").Value = "  +6.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.34"
$ws.Range("E37").Value = "  +1.39%  "

$ws.Range("E38").Value = "  +7.36%  "

$ws.Range("E39").Value = "  +5.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.27"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.92"
$ws.Range("E42").Value = "  +9.72%  "

$ws.Range("E43").Value = "  +7.96%  "

$ws.Range("E44").Value = "  +6.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.29"
$ws.Range("E45").Value = "  +23.10%  "

$ws.Range("E46").Value = "  +11.09%  "

$ws.Range("E47").Value = "  +7.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.01"
$ws.Range("E48").Value = "  +3.68%  "

$ws.Range("E49").Value = "  +16.65%  "

$ws.Range("E50").Value = "  +6.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.524"
$ws.Range("E51").Value = "  +6.51%  "
